# Append new rows (25-30) to tblStudyManagementTools sheet, matching the
# existing table's data pattern (A: numeric ID, B/C/D: text).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ A = 24; B = "18"; C = "Closures";    D = "Spatial Closures, of three types.  Marine protected areas, species risk based closures, and triggered closures based on effort" },
    @{ A = 25; B = "19"; C = "Catch Limit"; D = "four Landing regulations:  landing obligation, 5% discard limit, year-to-year quota transfer, both" },
    @{ A = 26; B = "20"; C = "Effort Limit"; D = "None, the current method is used in all simulations" },
    @{ A = 27; B = "21"; C = "Catch Limit"; D = "Status quo, 2x, Profit max, Broken stick control rule, Spatial broken stick control rule" },
    @{ A = 28; B = "21"; C = "Closures";    D = "Spatial closures.  Closure trigger by zone or species with 20% or 30% triggers" },
    @{ A = 29; B = "22"; C = "Catch Limit"; D = "8 approaches to setting P* buffer relative to F(lim)." }
)

$startRow = 25
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D

    $ws.Cells.Item($r, 1).Style = $ws.Cells.Item($r - 1, 1).Style
    $ws.Cells.Item($r, 2).Style = $ws.Cells.Item($r - 1, 2).Style
    $ws.Cells.Item($r, 3).Style = $ws.Cells.Item($r - 1, 3).Style
    $ws.Cells.Item($r, 4).Style = $ws.Cells.Item($r - 1, 4).Style
}
